# Update the "取得日時" (fetched-at) timestamp in column A for rows 2-10
# on the first worksheet ("ランサーズ") from 2026-01-25 18:27:13 to
# 2026-01-25 18:35:26, leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldTimestamp = "2026-01-25 18:27:13"
$newTimestamp = "2026-01-25 18:35:26"

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value() -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
